$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "d=6" distance column is inserted between the existing "d=5" (F)
# and "d=7" (old G) columns. The old G/H data columns shift right to H/I,
# and the new G column gets its own freshly computed values. Literal
# values are used throughout (rather than reading back existing cell
# values) to avoid any lossy read/round-trip of the stored doubles.

# ---- Header row (row 1) ----
$ws.Range("I1").Value = "d=10"
$ws.Range("H1").Value = "d=7"
$ws.Range("G1").Value = "d=6"

# Match the bold/border/centered header formatting used by the other
# header cells.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)

# ---- Row 2 (exponential) ----
$ws.Range("I2").Value = 96.05520876326425
$ws.Range("H2").Value = 98.80967168531303
$ws.Range("G2").Value = 95.83815861018448

# ---- Row 3 (mixture) ----
$ws.Range("I3").Value = 95.12136927754187
$ws.Range("H3").Value = 97.96646863542702
$ws.Range("G3").Value = 99.02524496274168

# ---- Row 4 (normal) ----
$ws.Range("I4").Value = 94.8521361022333
$ws.Range("H4").Value = 96.86018438497923
$ws.Range("G4").Value = 98.09409263363455

# ---- Row 5 (t-student) ----
$ws.Range("I5").Value = 97.08071733675607
$ws.Range("H5").Value = 98.60089019564199
$ws.Range("G5").Value = 98.33581724694211

# ---- Row 6 (uniform) ----
$ws.Range("I6").Value = 94.78077109621393
$ws.Range("H6").Value = 98.15428500623874
$ws.Range("G6").Value = 98.32153906875332
